# "Generate Report for Archive" - refresh localization-status report data.
#
# The underlying source rows keyed by the *.md guid are reordered
# (6ea86943 and 8cd68c7b move up, 5ca8c613 moves to the bottom) and the
# in-flight items (6ea86943 / 8cd68c7b) pick up a real "In Translation"
# status plus refreshed handoff metadata, while the cell-anchored
# hyperlinks keep pointing at their original fixed targets per cell
# address (B2/B3/B4/B5 on Overview, A2/A3/A4/A5 on the per-locale sheets)
# but need their display text (and the Overview row values) refreshed to
# match the new content in that cell.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Hyperlinks.Delete()

$ov.Range("A3").Value2 = "6ea86943-b7e0-4f8f-8f2a-891b6309fb6a.md"
$ov.Range("B3").Value2 = "e2e\6ea86943-b7e0-4f8f-8f2a-891b6309fb6a.md"
$ov.Range("C3").Value2 = ".md"
$ov.Range("D3").Value2 = ""
$ov.Range("E3").Value2 = "In Translation"
$ov.Range("F3").Value2 = "In Translation"
$ov.Range("G3").Value2 = "2016-08-25 12:42:45"

$ov.Range("A4").Value2 = "8cd68c7b-8843-4ec8-bb6a-1519158efe00.md"
$ov.Range("B4").Value2 = "e2e\8cd68c7b-8843-4ec8-bb6a-1519158efe00.md"
$ov.Range("C4").Value2 = ".md"
$ov.Range("D4").Value2 = ""
$ov.Range("E4").Value2 = "In Translation"
$ov.Range("F4").Value2 = "In Translation"
$ov.Range("G4").Value2 = "2016-08-25 12:42:45"

$ov.Range("A5").Value2 = "5ca8c613-518e-4f0d-bf89-3cdad89a4e79.md"
$ov.Range("B5").Value2 = "e2e\5ca8c613-518e-4f0d-bf89-3cdad89a4e79.md"
$ov.Range("C5").Value2 = ".md"
$ov.Range("D5").Value2 = ""
$ov.Range("E5").Value2 = "Ready for handoff"
$ov.Range("F5").Value2 = "Ready for handoff"
$ov.Range("G5").Value2 = "2016-08-25 12:41:32"

$ov.Hyperlinks.Add($ov.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f8c82938a2e1f926a623c832059cd149bae8af50/e2e/77229961-2680-49d2-952a-cce08c9f0f3f.md", [Type]::Missing, [Type]::Missing, "e2e\77229961-2680-49d2-952a-cce08c9f0f3f.md") | Out-Null
$ov.Hyperlinks.Add($ov.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/05d6eee429fd4b2e4940b3f1fb2f49c10429ed53/e2e/5ca8c613-518e-4f0d-bf89-3cdad89a4e79.md", [Type]::Missing, [Type]::Missing, "e2e\6ea86943-b7e0-4f8f-8f2a-891b6309fb6a.md") | Out-Null
$ov.Hyperlinks.Add($ov.Range("B4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ae96b5332522dfdba4240650b98d7dc1d647a914/e2e/6ea86943-b7e0-4f8f-8f2a-891b6309fb6a.md", [Type]::Missing, [Type]::Missing, "e2e\8cd68c7b-8843-4ec8-bb6a-1519158efe00.md") | Out-Null
$ov.Hyperlinks.Add($ov.Range("B5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ae96b5332522dfdba4240650b98d7dc1d647a914/e2e/8cd68c7b-8843-4ec8-bb6a-1519158efe00.md", [Type]::Missing, [Type]::Missing, "e2e\5ca8c613-518e-4f0d-bf89-3cdad89a4e79.md") | Out-Null

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Hyperlinks.Delete()

$zh.Range("A3").Value2 = "6ea86943-b7e0-4f8f-8f2a-891b6309fb6a.md"
$zh.Range("C3").Value2 = "In Translation"
$zh.Range("G3").Value2 = "6ea86943-b7e0-4f8f-8f2a-891b6309fb6a.f35575f34d2a966988c3460b78ce4689d5c30113.zh-cn.xlf"
$zh.Range("H3").Value2 = "2016-08-25 12:42:41"

$zh.Range("A4").Value2 = "8cd68c7b-8843-4ec8-bb6a-1519158efe00.md"
$zh.Range("C4").Value2 = "In Translation"
$zh.Range("G4").Value2 = "8cd68c7b-8843-4ec8-bb6a-1519158efe00.e78b1cc9c9ff3b089e1c7c5ab350849947932a89.zh-cn.xlf"
$zh.Range("H4").Value2 = "2016-08-25 12:42:41"

$zh.Range("A5").Value2 = "5ca8c613-518e-4f0d-bf89-3cdad89a4e79.md"
$zh.Range("C5").Value2 = "Ready for handoff"
$zh.Range("G5").Value2 = "5ca8c613-518e-4f0d-bf89-3cdad89a4e79.fe04a4cf668d341094941bcf5d58695972f0382c.zh-cn.xlf"
$zh.Range("H5").Value2 = "2016-08-25 12:41:27"

$zh.Hyperlinks.Add($zh.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f8c82938a2e1f926a623c832059cd149bae8af50/e2e/77229961-2680-49d2-952a-cce08c9f0f3f.md", [Type]::Missing, [Type]::Missing, "77229961-2680-49d2-952a-cce08c9f0f3f.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/a0585fee31c4c9e55085cb82f6c48e8f88eea701/e2e/77229961-2680-49d2-952a-cce08c9f0f3f.md", [Type]::Missing, [Type]::Missing, "77229961-2680-49d2-952a-cce08c9f0f3f.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/05d6eee429fd4b2e4940b3f1fb2f49c10429ed53/e2e/5ca8c613-518e-4f0d-bf89-3cdad89a4e79.md", [Type]::Missing, [Type]::Missing, "6ea86943-b7e0-4f8f-8f2a-891b6309fb6a.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ae96b5332522dfdba4240650b98d7dc1d647a914/e2e/6ea86943-b7e0-4f8f-8f2a-891b6309fb6a.md", [Type]::Missing, [Type]::Missing, "8cd68c7b-8843-4ec8-bb6a-1519158efe00.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ae96b5332522dfdba4240650b98d7dc1d647a914/e2e/8cd68c7b-8843-4ec8-bb6a-1519158efe00.md", [Type]::Missing, [Type]::Missing, "5ca8c613-518e-4f0d-bf89-3cdad89a4e79.md") | Out-Null

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Hyperlinks.Delete()

$de.Range("A3").Value2 = "6ea86943-b7e0-4f8f-8f2a-891b6309fb6a.md"
$de.Range("C3").Value2 = "In Translation"
$de.Range("G3").Value2 = "6ea86943-b7e0-4f8f-8f2a-891b6309fb6a.f35575f34d2a966988c3460b78ce4689d5c30113.de-de.xlf"
$de.Range("H3").Value2 = "2016-08-25 12:42:45"

$de.Range("A4").Value2 = "8cd68c7b-8843-4ec8-bb6a-1519158efe00.md"
$de.Range("C4").Value2 = "In Translation"
$de.Range("G4").Value2 = "8cd68c7b-8843-4ec8-bb6a-1519158efe00.e78b1cc9c9ff3b089e1c7c5ab350849947932a89.de-de.xlf"
$de.Range("H4").Value2 = "2016-08-25 12:42:45"

$de.Range("A5").Value2 = "5ca8c613-518e-4f0d-bf89-3cdad89a4e79.md"
$de.Range("C5").Value2 = "Ready for handoff"
$de.Range("G5").Value2 = "5ca8c613-518e-4f0d-bf89-3cdad89a4e79.fe04a4cf668d341094941bcf5d58695972f0382c.de-de.xlf"
$de.Range("H5").Value2 = "2016-08-25 12:41:32"

$de.Hyperlinks.Add($de.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f8c82938a2e1f926a623c832059cd149bae8af50/e2e/77229961-2680-49d2-952a-cce08c9f0f3f.md", [Type]::Missing, [Type]::Missing, "77229961-2680-49d2-952a-cce08c9f0f3f.md") | Out-Null
$de.Hyperlinks.Add($de.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/6a3daed74e14670b62aa64723d31f753b1ade2b6/e2e/77229961-2680-49d2-952a-cce08c9f0f3f.md", [Type]::Missing, [Type]::Missing, "77229961-2680-49d2-952a-cce08c9f0f3f.md") | Out-Null
$de.Hyperlinks.Add($de.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/05d6eee429fd4b2e4940b3f1fb2f49c10429ed53/e2e/5ca8c613-518e-4f0d-bf89-3cdad89a4e79.md", [Type]::Missing, [Type]::Missing, "6ea86943-b7e0-4f8f-8f2a-891b6309fb6a.md") | Out-Null
$de.Hyperlinks.Add($de.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ae96b5332522dfdba4240650b98d7dc1d647a914/e2e/6ea86943-b7e0-4f8f-8f2a-891b6309fb6a.md", [Type]::Missing, [Type]::Missing, "8cd68c7b-8843-4ec8-bb6a-1519158efe00.md") | Out-Null
$de.Hyperlinks.Add($de.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ae96b5332522dfdba4240650b98d7dc1d647a914/e2e/8cd68c7b-8843-4ec8-bb6a-1519158efe00.md", [Type]::Missing, [Type]::Missing, "5ca8c613-518e-4f0d-bf89-3cdad89a4e79.md") | Out-Null
